# Update FFXIV leve market-price data (H,I,J,K,L,M,N columns) per scheduled refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 133
$ws.Range("H133").Value = 14733.333
$ws.Range("J133").Value = 14733.333
$ws.Range("L133").Value = 14733.333
$ws.Range("N133").Value = -24853.333

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2647.3333
$ws.Range("I2").Value = 3693
$ws.Range("J2").Value = 1732.375
$ws.Range("K2").Value = 3693
$ws.Range("L2").Value = 1732.375
$ws.Range("M2").Value = -3580
$ws.Range("N2").Value = -1958.375
# Row 32
$ws.Range("H32").Value = 3862.9834
$ws.Range("I32").Value = 2302.2246
$ws.Range("J32").Value = 10815.454
$ws.Range("K32").Value = 2302.2246
$ws.Range("L32").Value = 10815.454
$ws.Range("M32").Value = -2015.2246
$ws.Range("N32").Value = -11389.454
# Row 61
$ws.Range("H61").Value = 3419.5518
$ws.Range("I61").Value = 2112.3076
$ws.Range("J61").Value = 4481.6875
$ws.Range("K61").Value = 2112.3076
$ws.Range("L61").Value = 4481.6875
$ws.Range("M61").Value = -1900.3076
$ws.Range("N61").Value = -4905.6875
# Row 97
$ws.Range("H97").Value = 20840078
$ws.Range("I97").Value = 27786418
$ws.Range("J97").Value = 1052.25
$ws.Range("K97").Value = 27786418
$ws.Range("L97").Value = 1052.25
$ws.Range("M97").Value = -27785922
$ws.Range("N97").Value = -2044.25
# Row 116
$ws.Range("H116").Value = 2647.3333
$ws.Range("I116").Value = 3693
$ws.Range("J116").Value = 1732.375
$ws.Range("K116").Value = 3693
$ws.Range("L116").Value = 1732.375
$ws.Range("M116").Value = -1399
$ws.Range("N116").Value = -6320.375
# Row 136
$ws.Range("H136").Value = 3419.5518
$ws.Range("I136").Value = 2112.3076
$ws.Range("J136").Value = 4481.6875
$ws.Range("K136").Value = 6336.9228
$ws.Range("L136").Value = 13445.0625
$ws.Range("M136").Value = -3786.9228
$ws.Range("N136").Value = -18545.0625
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2647.3333
$ws.Range("I3").Value = 3693
$ws.Range("J3").Value = 1732.375
$ws.Range("K3").Value = 3693
$ws.Range("L3").Value = 1732.375
$ws.Range("M3").Value = -3579
$ws.Range("N3").Value = -1960.375
# Row 82
$ws.Range("H82").Value = 16230.429
$ws.Range("I82").Value = 2609.4
$ws.Range("J82").Value = 50283
$ws.Range("K82").Value = 2609.4
$ws.Range("L82").Value = 50283
$ws.Range("M82").Value = -2226.4
$ws.Range("N82").Value = -51049
# Row 85
$ws.Range("H85").Value = 16230.429
$ws.Range("I85").Value = 2609.4
$ws.Range("J85").Value = 50283
$ws.Range("K85").Value = 2609.4
$ws.Range("L85").Value = 50283
$ws.Range("M85").Value = -1283.4
$ws.Range("N85").Value = -52935

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 783.3889
$ws.Range("I22").Value = 770.64703
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 770.64703
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -420.64703
$ws.Range("N22").Value = -1700

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 973.6667
$ws.Range("I5").Value = 717
$ws.Range("K5").Value = 2151
$ws.Range("M5").Value = -2039
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
# Row 26
$ws.Range("H26").Value = 240
$ws.Range("I26").Value = 240
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 720
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -432
$ws.Range("N26").ClearContents()
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
# Row 34
$ws.Range("H34").Value = 2914.1428
$ws.Range("J34").Value = 3299.8333
$ws.Range("L34").Value = 9899.499899999999
$ws.Range("N34").Value = -10067.4999
# Row 39
$ws.Range("H39").Value = 7502.6484
$ws.Range("I39").Value = 2666.6667
$ws.Range("J39").Value = 7929.353
$ws.Range("K39").Value = 8000.000100000001
$ws.Range("L39").Value = 23788.059
$ws.Range("M39").Value = -7706.000100000001
$ws.Range("N39").Value = -24376.059
# Row 40
$ws.Range("H40").Value = 251.33333
$ws.Range("I40").Value = 157.5
$ws.Range("J40").Value = 1002
$ws.Range("K40").Value = 630
$ws.Range("L40").Value = 4008
$ws.Range("M40").Value = -561
$ws.Range("N40").Value = -4146
# Row 46
$ws.Range("H46").Value = 775
$ws.Range("I46").Value = 550
$ws.Range("K46").Value = 1650
$ws.Range("M46").Value = -1559
# Row 122
$ws.Range("H122").Value = 629.6
$ws.Range("I122").Value = 283.83334
$ws.Range("J122").Value = 1148.25
$ws.Range("K122").Value = 2554.50006
$ws.Range("L122").Value = 10334.25
$ws.Range("M122").Value = -104.5000600000003
$ws.Range("N122").Value = -15234.25
# Row 135
$ws.Range("H135").Value = 973.6667
$ws.Range("I135").Value = 717
$ws.Range("K135").Value = 6453
$ws.Range("M135").Value = -3918

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6543
$ws.Range("I70").Value = 6890.375
$ws.Range("J70").Value = 5616.6665
$ws.Range("K70").Value = 6890.375
$ws.Range("L70").Value = 5616.6665
$ws.Range("M70").Value = -6620.375
$ws.Range("N70").Value = -6156.6665
# Row 73
$ws.Range("H73").Value = 6543
$ws.Range("I73").Value = 6890.375
$ws.Range("J73").Value = 5616.6665
$ws.Range("K73").Value = 6890.375
$ws.Range("L73").Value = 5616.6665
$ws.Range("M73").Value = -5954.375
$ws.Range("N73").Value = -7488.6665
# Row 80
$ws.Range("H80").Value = 2233.75
$ws.Range("I80").Value = 1977
$ws.Range("J80").Value = 2661.6667
$ws.Range("K80").Value = 1977
$ws.Range("L80").Value = 2661.6667
$ws.Range("M80").Value = -979
$ws.Range("N80").Value = -4657.6667
# Row 83
$ws.Range("H83").Value = 2233.75
$ws.Range("I83").Value = 1977
$ws.Range("J83").Value = 2661.6667
$ws.Range("K83").Value = 9885
$ws.Range("L83").Value = 13308.3335
$ws.Range("M83").Value = -4893
$ws.Range("N83").Value = -23292.3335
# Row 123
$ws.Range("H123").Value = 10638.5
$ws.Range("J123").Value = 10638.5
$ws.Range("L123").Value = 10638.5
$ws.Range("N123").Value = -15538.5
# Row 137
$ws.Range("H137").Value = 34597.5
$ws.Range("J137").Value = 34597.5
$ws.Range("L137").Value = 34597.5
$ws.Range("N137").Value = -44797.5
# Row 138
$ws.Range("H138").Value = 31533.334
$ws.Range("J138").Value = 31533.334
$ws.Range("L138").Value = 31533.334
$ws.Range("N138").Value = -41813.334
# Row 139
$ws.Range("H139").Value = 30000
$ws.Range("J139").Value = 30000
$ws.Range("L139").Value = 30000
$ws.Range("N139").Value = -40280

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 380.15384
$ws.Range("I55").Value = 274
$ws.Range("J55").Value = 734
$ws.Range("K55").Value = 274
$ws.Range("L55").Value = 734
$ws.Range("M55").Value = -101
$ws.Range("N55").Value = -1080
# Row 56
$ws.Range("H56").Value = 3050.1428
$ws.Range("J56").Value = 3075
$ws.Range("L56").Value = 3075
$ws.Range("N56").Value = -4457
# Row 93
$ws.Range("H93").Value = 976
$ws.Range("J93").Value = 952
$ws.Range("L93").Value = 952
$ws.Range("N93").Value = -3448
